# Edit the "area_pop_sum" sheet: drop the "Density" column (C), fix casing of
# row labels, and re-lay the density value as a fourth row instead of a third
# column.

$wb = $excel.ActiveWorkbook
$origActive = $wb.ActiveSheet
$ws = $wb.Worksheets.Item("area_pop_sum")

# Capture the density value (currently duplicated at B-less column C) before
# we start clearing cells.
$densityValue = $ws.Range("C2").Value2

# Remove column C entirely (drops C1 "Density" header, C2 and C3 values) and
# shift remaining columns left.
$ws.Range("C:C").Delete()

# Relabel the row-3 index from "Population" to lowercase "population".
$ws.Range("A3").Value = "population"

# Append the density figure as a new row instead of a third column.
$ws.Range("A4").Value = "density"
$ws.Range("B4").Value = $densityValue

$ws.Range("A1").Select()

# Restore whichever sheet/tab was active before this edit ran.
$origActive.Activate()
